# Segundo happy test validado
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Happy Path")

# Fill in row 3 (test case #2 - "Agregar un item a el carrito")
$ws.Range("B3").Value = "Agregar un item a el carrito"
$ws.Range("C3").Value = "Ir a la página del sitio"
$ws.Range("E3").Value = "1. Hacer hover en la tarjeta de algún item`r`n2. Hacer click en ""Add to cart"""
$ws.Range("F3").Value = "Aparece un modal con la notificación de que el item fue agregado al carrito."
$ws.Range("G3").Value = $ws.Range("G2").Value
$ws.Range("H3").Value = $ws.Range("H2").Value

# Formatting to match the filled-in row
$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("B3").WrapText = $false

$ws.Range("C3").HorizontalAlignment = -4108
$ws.Range("C3").VerticalAlignment = -4108
$ws.Range("C3").WrapText = $false

$ws.Range("E3").HorizontalAlignment = -4131
$ws.Range("E3").VerticalAlignment = -4108
$ws.Range("E3").WrapText = $true

$ws.Columns.Item(2).ColumnWidth = 25.7109375
$ws.Columns.Item(2).HorizontalAlignment = -4108
$ws.Columns.Item(2).VerticalAlignment = -4108

# Add a new row 6 for test case ID 5 (blank, following the same pattern)
$ws.Range("A6").Value = 5
$ws.Range("B6").HorizontalAlignment = -4108
$ws.Range("B6").VerticalAlignment = -4108
$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("D6").VerticalAlignment = -4108
$ws.Range("E6").HorizontalAlignment = -4131
$ws.Range("E6").VerticalAlignment = -4108

# Rows 4 and 5 B/D/E alignment updates
$ws.Range("B4").HorizontalAlignment = -4108
$ws.Range("B4").VerticalAlignment = -4108
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D4").VerticalAlignment = -4108
$ws.Range("E4").HorizontalAlignment = -4131
$ws.Range("E4").VerticalAlignment = -4108

$ws.Range("B5").HorizontalAlignment = -4108
$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("D5").HorizontalAlignment = -4108
$ws.Range("D5").VerticalAlignment = -4108
$ws.Range("E5").HorizontalAlignment = -4131
$ws.Range("E5").VerticalAlignment = -4108

$ws.Range("B3").Select()
